# Week 3 materials: insert a new "Variables" section-divider slide at
# position 4 (Title Only layout) and refresh the auto date-figure
# placeholder text (9.03.2025 -> 12.03.2025) across the slide master
# and every slide layout.

$p = $ppt.ActivePresentation

# --- refresh the cached "today" date text on the master ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "12.03.2025"
    }
}

# --- same refresh on every slide layout ---
for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "12.03.2025"
        }
    }
}

# --- insert the new "Variables" divider slide at position 4 ---
# ppLayoutTitleOnly = 11
$newSlide = $p.Slides.Add(4, 11)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Variables"
